$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A11").Value = "Hazan"
$ws.Range("B11").Value = "Elad"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 2

$ws.Range("A12").Value = "Kingma"
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 2

$ws.Range("A13").Value = "Ba"
$ws.Range("B13").Value = "Jimmy"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 2

$ws.Range("B12").Select()
